# Rename sheet "strategy_id-5008" to "strategy_id-5007" and add a new
# sheet "strategy_id-5009" (an identical copy of it) right after it.
$wb = $excel.ActiveWorkbook

$src = $wb.Worksheets.Item("strategy_id-5008")

# Duplicate the existing sheet, placing the copy right after the source.
# The copy will pick up the exact same data/styles as the source sheet.
$src.Copy($null, $src)

# Rename the original sheet to its new name.
$src.Name = "strategy_id-5007"

# The freshly created copy is placed immediately after $src and is named
# "strategy_id-5008 (2)" by default - rename it to the new sheet name.
$newSheet = $wb.Worksheets.Item("strategy_id-5008 (2)")
$newSheet.Name = "strategy_id-5009"

# Keep the first sheet as the active tab, matching the original workbook.
$wb.Worksheets.Item(1).Activate()
